$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add 12 new rows of training-log data (rows 830-841), mirroring the
# formatting of the existing rows 822 (empty "Localisation douleur")
# and 829 (populated "Localisation douleur").
# ------------------------------------------------------------------

# Template rows already present on the sheet, used purely to carry the
# correct cell styles (s="3" dates, s="1" values, s="2" empty G cells)
# onto the freshly created rows.
$emptyGTemplate = "A822:I822"   # row with an empty "G" (Localisation douleur) cell
$filledGTemplate = "G829"       # single cell with a populated/styled "G" cell

$rows = @(830,831,832,833,834,835,836,837,838,839,840,841)
foreach ($r in $rows) {
    $ws.Range($emptyGTemplate).Copy($ws.Range("A$r`:I$r"))
}

# NOTE: named PowerShell parameters are not reliably bound by this
# runtime, so Set-LogRow uses plain positional parameters instead.
function Set-LogRow {
    param($Row, $Date, $Name, $Volume, $Intensite, $Fatigue, $Douleur, $Localisation, $Plaisir)

    $ws.Range("A$Row").Value = $Date
    $ws.Range("B$Row").Value = $Name
    $ws.Range("C$Row").Value = $Volume
    $ws.Range("D$Row").Value = $Intensite
    $ws.Range("E$Row").Value = $Fatigue
    $ws.Range("F$Row").Value = $Douleur
    if ([string]::IsNullOrEmpty($Localisation)) {
        $ws.Range("G$Row").ClearContents()
    } else {
        $ws.Range($filledGTemplate).Copy($ws.Range("G$Row"))
        $ws.Range("G$Row").Value = $Localisation
    }
    $ws.Range("H$Row").Value = $Plaisir
}

Set-LogRow 830 46064 "Omar Benyounes"  70 7 5 0 ""                            8
Set-LogRow 831 46064 "Yoan Zouma"      70 5 8 5 "Fesse"                       3
Set-LogRow 832 46064 "Kamal Bafounta"  70 5 5 0 ""                            4
Set-LogRow 833 46064 "Naim Ighbane"    70 7 8 7 "Courbature"                  5
Set-LogRow 834 46064 "Jeremie Laurent" 70 7 7 1 "Adducteur "           7
Set-LogRow 835 46064 "Yoann Martelat"  70 6 7 6 "Genou"                       2
Set-LogRow 836 46064 "Levy Ndoutoume"  70 6 7 0 ""                            7
Set-LogRow 837 46064 "Ilan Ihaddadene" 70 6 7 0 ""                            9
Set-LogRow 838 46064 "Karahali Souaré" 70 5 7 6 "Cheville"                    5
Set-LogRow 839 46064 "Theo Owono"      70 5 4 3 "Fesse"                       2
Set-LogRow 840 46064 "Romain Thunet"   70 5 5 3 "Synthétique "         7
Set-LogRow 841 46064 "Nathanael Beta"  70 6 6 4 "Dos et courbature "          7

# ------------------------------------------------------------------
# Charge (column I) = Volume * Intensite, filled down as two shared
# formula blocks matching the break introduced by the source workbook.
# ------------------------------------------------------------------
$ws.Range("I830:I835").Formula = "=C830*D830"
$ws.Range("I836:I841").Formula = "=C836*D836"

# ------------------------------------------------------------------
# Restore the scroll position / selection recorded in the workbook.
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 811
$win.ScrollColumn = 1
$ws.Range("K834").Select()
